# Update the "Förändrad" (changed) date column C for all data rows (2-54)
# from 2023-09-16 (45185) to 2023-10-05 (45204).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 54
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45185) {
        $cell.Value2 = 45204
    }
}
